# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Updates the computed "K" values for rows 2-28 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 3
    4  = 9
    5  = 6
    6  = 8
    7  = 5
    8  = 4
    9  = 7
    10 = 4
    11 = 6
    12 = 3
    13 = 4
    14 = 6
    15 = 3
    16 = 4
    17 = 6
    18 = 5
    19 = 3
    20 = 6
    21 = 5
    22 = 6
    23 = 4
    24 = 2
    25 = 7
    26 = 5
    27 = 1
    28 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
